$wb = $excel.ActiveWorkbook

function Convert-CellToText($cell) {
    $val = $cell.Value()
    $formatted = $val.ToString("#,##0")
    $cell.NumberFormat = "@"
    $cell.Value = $formatted
}

function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# --- Sheet "Overall": A2 (count of 990 filers) becomes text "3,468" ---
$wsOverall = $wb.Worksheets.Item("Overall")
Convert-CellToText $wsOverall.Cells.Item(2, 1)

# --- Sheet "County": B2:B99 (filer counts per county) become text ---
$wsCounty = $wb.Worksheets.Item("County")
for ($r = 2; $r -le 99; $r++) {
    Convert-CellToText $wsCounty.Cells.Item($r, 2)
}

# --- Sheet "County": append the missing "Total" row (row 100) ---
Set-TextValue $wsCounty.Cells.Item(100, 1) "Total"
Set-TextValue $wsCounty.Cells.Item(100, 2) "3,468"
Set-TextValue $wsCounty.Cells.Item(100, 3) '$9,986,549,815'
Set-TextValue $wsCounty.Cells.Item(100, 4) "8.22%"
Set-TextValue $wsCounty.Cells.Item(100, 5) "-13.23%"
Set-TextValue $wsCounty.Cells.Item(100, 6) "67.30%"

# --- Sheet "Congressional District": B2:B19 become text ---
$wsCd = $wb.Worksheets.Item("Congressional District")
for ($r = 2; $r -le 19; $r++) {
    Convert-CellToText $wsCd.Cells.Item($r, 2)
}

# --- Sheet "Size": B2:B8 become text ---
$wsSize = $wb.Worksheets.Item("Size")
for ($r = 2; $r -le 8; $r++) {
    Convert-CellToText $wsSize.Cells.Item($r, 2)
}

# --- Sheet "Subsector": B2:B14 become text ---
$wsSubsector = $wb.Worksheets.Item("Subsector")
for ($r = 2; $r -le 14; $r++) {
    Convert-CellToText $wsSubsector.Cells.Item($r, 2)
}

Write-Host "done"
